$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Slide 6: change the table's style (tableStyleId) from the custom
#    "Table_0" style to the built-in style {A5A373D5-B56A-4664-87CC-A54C0006EFDF}.
#    The table is the 2nd shape on slide 6 (the 1st shape is the heading text).
# ---------------------------------------------------------------------------
$s = $p.Slides.Item(6)
$tbl = $s.Shapes.Item(2).Table
$tbl.ApplyStyle("{A5A373D5-B56A-4664-87CC-A54C0006EFDF}")

# ---------------------------------------------------------------------------
# 2) Re-colour the presentation's theme from "Integral" to the stock
#    "Office Theme" colour palette (dk1/lt1/dk2/lt2/accent1-6/hlink/folHlink).
#    RGB values are supplied as OLE colours (0x00BBGGRR).
# ---------------------------------------------------------------------------
$theme = $p.Designs.Item(1).SlideMaster.Theme
$colors = $theme.ThemeColorScheme

$colors.Item(1).RGB  = 0x000000   # dk1
$colors.Item(2).RGB  = 0xFFFFFF   # lt1
$colors.Item(3).RGB  = 0x6A5444   # dk2      (44546A)
$colors.Item(4).RGB  = 0xE6E6E7   # lt2      (E7E6E6)
$colors.Item(5).RGB  = 0xD59B5B   # accent1  (5B9BD5)
$colors.Item(6).RGB  = 0x317DED   # accent2  (ED7D31)
$colors.Item(7).RGB  = 0xA5A5A5   # accent3  (A5A5A5)
$colors.Item(8).RGB  = 0x00C0FF   # accent4  (FFC000)
$colors.Item(9).RGB  = 0xC47244   # accent5  (4472C4)
$colors.Item(10).RGB = 0x47AD70   # accent6  (70AD47)
$colors.Item(11).RGB = 0xC16305   # hlink    (0563C1)
$colors.Item(12).RGB = 0x724F95   # folHlink (954F72)
